$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.248.42"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "3.226.00"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'606.79"
$ws.Range("E5").Value = "  +5.44%  "
$ws.Range("D6").Value = "'154.24"
$ws.Range("E6").Value = "  +2.68%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.222.34"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").Value = "'0.537"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "'6.19"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "'0.515"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "'39.02"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").Value = "3.753.68"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").Value = "66.225.47"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "'7.46"
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("D18").Value = "3.231.62"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "'515.52"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  +6.36%  "
$ws.Range("D22").Value = "'0.743"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").Value = "'15.31"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'8.04"
$ws.Range("E24").Value = "  +2.82%  "
$ws.Range("D25").Value = "'85.73"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'3.05"
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'9.28"
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("D30").Value = "'2.93"
$ws.Range("E30").Value = "  +5.79%  "
$ws.Range("D31").Value = "'6.86"
$ws.Range("E31").Value = "  +9.03%  "
$ws.Range("D32").Value = "'28.39"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "'6.69"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "'55.60"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'0.0932"
$ws.Range("E37").Value = "  +4.48%  "
$ws.Range("D38").Value = "'492.26"
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("D39").Value = "'0.0425"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").Value = "'3.04"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").Value = "'8.91"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.296"
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.040.22"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").Value = "0.0₃0648"
$ws.Range("E45").Value = "  +8.38%  "
$ws.Range("D46").Value = "'2.47"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").Value = "'29.32"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D49").Value = "'0.117"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").Value = "'2.35"
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "'119.98"
$ws.Range("E51").Value = "  -1.04%  "
